$wb = $excel.ActiveWorkbook

# The F column ("想去人数" / "wanted to go" count) was refreshed for both the
# "展览" sheet and the mirrored "全部类型" sheet. Apply the same updates to
# each worksheet that contains this data.
$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 1075
    7  = 2369
    11 = 1167
    14 = 3
    15 = 1056
    17 = 305
    18 = 14
    22 = 79
    24 = 19
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
